$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update predictor names in column C to use ln(...) notation with bracketed units,
# per the shared-strings changes in the diff.
$ws.Range("C2").Value  = "ln(GDP [dollars per capita])"
$ws.Range("C3").Value  = "ln(Tourism - Inbound [per capita])"
$ws.Range("C4").Value  = "ln(ProMed Mentions [per capita])"
$ws.Range("C8").Value  = "ln(Publication Bias Index [per capita])"
$ws.Range("C9").Value  = "ln(AB Exports [dollars per capita])"
$ws.Range("C10").Value = "ln(Livestock AB Consumption [per PCU])"
$ws.Range("C11").Value = "ln(Migrant Population [per capita])"
$ws.Range("C12").Value = "ln(Livestock Population [PCU])"
$ws.Range("C13").Value = "ln(ProMed Mentions [per capita])"
$ws.Range("C15").Value = "ln(Publication Bias Index [per capita])"
$ws.Range("C16").Value = "ln(Population)"
$ws.Range("C17").Value = "ln(GDP [dollars per capita])"
